$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 247
$ws1.Range("F4").Value = 810
$ws1.Range("F5").Value = 518

# Sheet "全部类型": update "想去人数" (F column) values
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F2").Value = 247
$ws2.Range("F4").Value = 810
$ws2.Range("F6").Value = 518
